# aggiornamento fino a 1/09/2021
# Append new daily rows (358-366) to Sheet1, mirroring the existing
# layout: col A = date serial (formatted/bordered/centered like the
# preceding rows), col B = nuovi positivi, col C = somma mobile 7gg.,
# col D = somma mobile 7gg. per 100mila abitanti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(44432, 2, 7, 116.4531691898187),
    @(44433, 0, 7, 116.4531691898187),
    @(44434, 0, 6, 99.81700216270171),
    @(44435, 7, 10, 166.3616702711695),
    @(44436, 0, 10, 166.3616702711695),
    @(44437, 1, 10, 166.3616702711695),
    @(44438, 2, 12, 199.6340043254034),
    @(44439, 3, 13, 216.2701713525204),
    @(44440, 0, 13, 216.2701713525204)
)

$startRow = 358
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $data[0]
    $cellA.Font.Bold = $true
    $cellA.Borders.LineStyle = 1
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
